$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44413
$ws.Range("J2").Value = 32
$ws.Range("K2").Value = 9000
$ws.Range("M2").Value = 9438
$ws.Range("P2").Value = 3146
$ws.Range("D3").Value = 44203
$ws.Range("J3").Value = 109
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 3333
$ws.Range("D4").Value = 44414
$ws.Range("J4").Value = 29
$ws.Range("K4").Value = 9000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 9483
$ws.Range("P4").Value = 3161
$ws.Range("D5").Value = 44217
$ws.Range("J5").Value = 35
$ws.Range("D6").Value = 44434
$ws.Range("J6").Value = 32
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 9562
$ws.Range("O6").Value = 'Provincia de Santiago'
$ws.Range("P6").Value = 3187
$ws.Range("D7").Value = 44446
$ws.Range("J7").Value = 32
$ws.Range("K7").Value = 9000
$ws.Range("M7").Value = 9469
$ws.Range("P7").Value = 3156
$ws.Range("D8").Value = 44161
$ws.Range("J8").Value = 37
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 11595
$ws.Range("P8").Value = 3865
$ws.Range("D9").Value = 44257
$ws.Range("J9").Value = 42
$ws.Range("D10").Value = 44504
$ws.Range("J10").Value = 35
$ws.Range("K10").Value = 7500
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 7729
$ws.Range("P10").Value = 2576
$ws.Range("D11").Value = 44271
$ws.Range("J11").Value = 36
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 3333
$ws.Range("D12").Value = 44260
$ws.Range("J12").Value = 33
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("P12").Value = 3333
$ws.Range("D13").Value = 44264
$ws.Range("K13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("O13").Value = 'Región Metropolitana'
$ws.Range("P13").Value = 3333
$ws.Range("D14").Value = 44516
$ws.Range("J14").Value = 34
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 9000
$ws.Range("M14").Value = 8441
$ws.Range("P14").Value = 2814
$ws.Range("D15").Value = 44488
$ws.Range("J15").Value = 34
$ws.Range("K15").Value = 8500
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 8735
$ws.Range("P15").Value = 2912
$ws.Range("D16").Value = 44425
$ws.Range("K16").Value = 9000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 9438
$ws.Range("P16").Value = 3146
$ws.Range("D17").Value = 44187
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 10000
$ws.Range("P17").Value = 3333
$ws.Range("D18").Value = 44523
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 8500
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8729
$ws.Range("P18").Value = 2910
$ws.Range("D19").Value = 44343
$ws.Range("J19").Value = 37
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 8500
$ws.Range("M19").Value = 8203
$ws.Range("P19").Value = 2734
$ws.Range("D20").Value = 44467
$ws.Range("J20").Value = 33
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9515
$ws.Range("P20").Value = 3172
$ws.Range("D21").Value = 44448
$ws.Range("J21").Value = 32
$ws.Range("K21").Value = 8500
$ws.Range("M21").Value = 8734
$ws.Range("P21").Value = 2911
$ws.Range("D22").Value = 44334
$ws.Range("J22").Value = 39
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 8410
$ws.Range("P22").Value = 2803
$ws.Range("D23").Value = 44294
$ws.Range("J23").Value = 38
$ws.Range("D24").Value = 44175
$ws.Range("J24").Value = 41
$ws.Range("K24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 10000
$ws.Range("P24").Value = 3333
$ws.Range("D26").Value = 44330
$ws.Range("J26").Value = 45
$ws.Range("K26").Value = 8500
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = 8744
$ws.Range("P26").Value = 2915
$ws.Range("D27").Value = 44222
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = 10000
$ws.Range("M27").Value = 10000
$ws.Range("P27").Value = 3333
$ws.Range("D28").Value = 44238
$ws.Range("J28").Value = 35
$ws.Range("D29").Value = 44364
$ws.Range("J29").Value = 35
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 7457
$ws.Range("P29").Value = 2486
$ws.Range("D30").Value = 44351
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 7405
$ws.Range("P30").Value = 2468
$ws.Range("D31").Value = 44329
$ws.Range("J31").Value = 43
$ws.Range("K31").Value = 7500
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = 7733
$ws.Range("P31").Value = 2578
$ws.Range("D32").Value = 44266
$ws.Range("J32").Value = 43
$ws.Range("D33").Value = 44302
$ws.Range("J33").Value = 44
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 10000
$ws.Range("O33").Value = 'Región Metropolitana'
$ws.Range("P33").Value = 3333
$ws.Range("D34").Value = 44455
$ws.Range("J34").Value = 28
$ws.Range("K34").Value = 8500
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = 8732
$ws.Range("P34").Value = 2911
$ws.Range("D35").Value = 44369
$ws.Range("J35").Value = 35
$ws.Range("K35").Value = 7000
$ws.Range("L35").Value = 7500
$ws.Range("M35").Value = 7229
$ws.Range("P35").Value = 2410
$ws.Range("D36").Value = 44483
$ws.Range("J36").Value = 33
$ws.Range("K36").Value = 8500
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = 8727
$ws.Range("P36").Value = 2909
$ws.Range("D37").Value = 44524
$ws.Range("J37").Value = 34
$ws.Range("K37").Value = 8500
$ws.Range("M37").Value = 8721
$ws.Range("P37").Value = 2907
$ws.Range("D38").Value = 44460
$ws.Range("J38").Value = 29
$ws.Range("K38").Value = 9000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = 9483
$ws.Range("P38").Value = 3161
$ws.Range("D39").Value = 44442
$ws.Range("J39").Value = 32
$ws.Range("K39").Value = 9000
$ws.Range("M39").Value = 9562
$ws.Range("P39").Value = 3187
$ws.Range("D40").Value = 44320
$ws.Range("J40").Value = 42
$ws.Range("K40").Value = 9000
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = 9000
$ws.Range("P40").Value = 3000
$ws.Range("D41").Value = 44385
$ws.Range("J41").Value = 32
$ws.Range("K41").Value = 8500
$ws.Range("L41").Value = 9000
$ws.Range("M41").Value = 8719
$ws.Range("P41").Value = 2906
$ws.Range("D42").Value = 44194
$ws.Range("J42").Value = 45
$ws.Range("D43").Value = 44365
$ws.Range("J43").Value = 32
$ws.Range("K43").Value = 7500
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = 7734
$ws.Range("P43").Value = 2578
$ws.Range("D44").Value = 44495
$ws.Range("J44").Value = 35
$ws.Range("M44").Value = 8457
$ws.Range("P44").Value = 2819
$ws.Range("D45").Value = 44348
$ws.Range("J45").Value = 42
$ws.Range("M45").Value = 8405
$ws.Range("P45").Value = 2802
$ws.Range("D46").Value = 44511
$ws.Range("J46").Value = 35
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 9000
$ws.Range("M46").Value = 8457
$ws.Range("P46").Value = 2819
$ws.Range("D47").Value = 44309
$ws.Range("J47").Value = 39
$ws.Range("L47").Value = 9000
$ws.Range("M47").Value = 9000
$ws.Range("O47").Value = 'Región Metropolitana'
$ws.Range("P47").Value = 3000
$ws.Range("D48").Value = 44477
$ws.Range("J48").Value = 35
$ws.Range("K48").Value = 8500
$ws.Range("M48").Value = 8729
$ws.Range("P48").Value = 2910
$ws.Range("D49").Value = 44323
$ws.Range("J49").Value = 45
$ws.Range("K49").Value = 7500
$ws.Range("M49").Value = 7778
$ws.Range("P49").Value = 2593
$ws.Range("D50").Value = 44420
$ws.Range("J50").Value = 33
$ws.Range("K50").Value = 9000
$ws.Range("M50").Value = 9455
$ws.Range("O50").Value = 'Región Metropolitana'
$ws.Range("P50").Value = 3152
$ws.Range("D51").Value = 44328
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = 8421
$ws.Range("P51").Value = 2807
$ws.Range("D52").Value = 44295
$ws.Range("J52").Value = 44
$ws.Range("K52").Value = 9000
$ws.Range("M52").Value = 9568
$ws.Range("P52").Value = 3189
$ws.Range("D53").Value = 44246
$ws.Range("J53").Value = 37
$ws.Range("K53").Value = 10000
$ws.Range("M53").Value = 10000
$ws.Range("O53").Value = 'Región Metropolitana'
$ws.Range("P53").Value = 3333
$ws.Range("D54").Value = 44225
$ws.Range("J54").Value = 32
$ws.Range("K54").Value = 10000
$ws.Range("M54").Value = 10000
$ws.Range("P54").Value = 3333
$ws.Range("D55").Value = 44186
$ws.Range("J55").Value = 70
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = 10000
$ws.Range("P55").Value = 3333
$ws.Range("D56").Value = 44362
$ws.Range("J56").Value = 41
$ws.Range("K56").Value = 7500
$ws.Range("L56").Value = 8000
$ws.Range("M56").Value = 7720
$ws.Range("P56").Value = 2573
$ws.Range("D57").Value = 44376
$ws.Range("J57").Value = 32
$ws.Range("K57").Value = 9000
$ws.Range("L57").Value = 10000
$ws.Range("M57").Value = 9594
$ws.Range("P57").Value = 3198
$ws.Range("D58").Value = 44327
$ws.Range("J58").Value = 44
$ws.Range("K58").Value = 8000
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = 8455
$ws.Range("P58").Value = 2818
$ws.Range("D59").Value = 44308
$ws.Range("J59").Value = 37
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 10000
$ws.Range("M59").Value = 10000
$ws.Range("P59").Value = 3333
$ws.Range("D60").Value = 44411
$ws.Range("J60").Value = 32
$ws.Range("K60").Value = 9000
$ws.Range("M60").Value = 9469
$ws.Range("P60").Value = 3156
$ws.Range("D61").Value = 44176
$ws.Range("J61").Value = 35
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = 10000
$ws.Range("P61").Value = 3333
$ws.Range("D62").Value = 44530
$ws.Range("J62").Value = 34
$ws.Range("K62").Value = 9000
$ws.Range("M62").Value = 9471
$ws.Range("O62").Value = 'Provincia de Santiago'
$ws.Range("P62").Value = 3157
$ws.Range("D63").Value = 44336
$ws.Range("J63").Value = 43
$ws.Range("K63").Value = 8000
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = 8419
$ws.Range("P63").Value = 2806
$ws.Range("D64").Value = 44292
$ws.Range("J64").Value = 39
$ws.Range("K64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = 10000
$ws.Range("P64").Value = 3333
$ws.Range("D66").Value = 44245
$ws.Range("J66").Value = 37
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = 10000
$ws.Range("P66").Value = 3333
$ws.Range("D67").Value = 44243
$ws.Range("J67").Value = 28
$ws.Range("D68").Value = 44250
$ws.Range("J68").Value = 37
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = 10000
$ws.Range("O68").Value = 'Provincia de Santiago'
$ws.Range("P68").Value = 3333
$ws.Range("D69").Value = 44476
$ws.Range("J69").Value = 34
$ws.Range("K69").Value = 8500
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = 8735
$ws.Range("P69").Value = 2912
$ws.Range("D70").Value = 44236
$ws.Range("J70").Value = 36
$ws.Range("D71").Value = 44196
$ws.Range("J71").Value = 33
$ws.Range("D72").Value = 44306
$ws.Range("J72").Value = 42
$ws.Range("K72").Value = 10000
$ws.Range("M72").Value = 10000
$ws.Range("P72").Value = 3333
$ws.Range("D73").Value = 44188
$ws.Range("J73").Value = 70
$ws.Range("D74").Value = 44397
$ws.Range("J74").Value = 80
$ws.Range("K74").Value = 11000
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = 11375
$ws.Range("O74").Value = 'Provincia de Santiago'
$ws.Range("P74").Value = 3792
$ws.Range("D75").Value = 44285
$ws.Range("J75").Value = 37
$ws.Range("D76").Value = 44371
$ws.Range("J76").Value = 35
$ws.Range("K76").Value = 7000
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = 7229
$ws.Range("P76").Value = 2410
$ws.Range("D77").Value = 44390
$ws.Range("J77").Value = 28
$ws.Range("K77").Value = 8500
$ws.Range("M77").Value = 8732
$ws.Range("P77").Value = 2911
$ws.Range("D78").Value = 44159
$ws.Range("J78").Value = 34
$ws.Range("K78").Value = 12000
$ws.Range("L78").Value = 12000
$ws.Range("M78").Value = 12000
$ws.Range("P78").Value = 4000
$ws.Range("D79").Value = 44355
$ws.Range("J79").Value = 38
$ws.Range("K79").Value = 7500
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = 7724
$ws.Range("P79").Value = 2575
$ws.Range("D80").Value = 44406
$ws.Range("J80").Value = 31
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = 9452
$ws.Range("P80").Value = 3151
$ws.Range("D81").Value = 44384
$ws.Range("J81").Value = 30
$ws.Range("K81").Value = 9000
$ws.Range("M81").Value = 9000
$ws.Range("P81").Value = 3000
$ws.Range("D82").Value = 44432
$ws.Range("J82").Value = 31
$ws.Range("M82").Value = 9452
$ws.Range("P82").Value = 3151
$ws.Range("D83").Value = 44201
$ws.Range("J83").Value = 36
$ws.Range("D84").Value = 44208
$ws.Range("J84").Value = 88
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = 10000
$ws.Range("P84").Value = 3333
$ws.Range("D85").Value = 44210
$ws.Range("J85").Value = 70
$ws.Range("D86").Value = 44315
$ws.Range("J86").Value = 43
$ws.Range("K86").Value = 8000
$ws.Range("M86").Value = 8581
$ws.Range("P86").Value = 2860
$ws.Range("D87").Value = 44519
$ws.Range("J87").Value = 33
$ws.Range("K87").Value = 8000
$ws.Range("L87").Value = 9000
$ws.Range("M87").Value = 8455
$ws.Range("P87").Value = 2818
$ws.Range("D88").Value = 44462
$ws.Range("J88").Value = 32
$ws.Range("K88").Value = 9000
$ws.Range("M88").Value = 9531
$ws.Range("P88").Value = 3177
$ws.Range("D89").Value = 44168
$ws.Range("J89").Value = 39
$ws.Range("K89").Value = 12000
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = 12000
$ws.Range("P89").Value = 4000
$ws.Range("D90").Value = 44313
$ws.Range("J90").Value = 37
$ws.Range("L90").Value = 9000
$ws.Range("M90").Value = 9000
$ws.Range("P90").Value = 3000
$ws.Range("D91").Value = 44435
$ws.Range("J91").Value = 32
$ws.Range("M91").Value = 9562
$ws.Range("O91").Value = 'Provincia de Santiago'
$ws.Range("P91").Value = 3187
$ws.Range("D92").Value = 44435
$ws.Range("J92").Value = 63
$ws.Range("K92").Value = 9000
$ws.Range("M92").Value = 9460
$ws.Range("P92").Value = 3153
$ws.Range("D93").Value = 44322
$ws.Range("J93").Value = 42
$ws.Range("M93").Value = 8476
$ws.Range("P93").Value = 2825
$ws.Range("D94").Value = 44231
$ws.Range("J94").Value = 31
$ws.Range("K94").Value = 10000
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = 10000
$ws.Range("P94").Value = 3333
$ws.Range("D95").Value = 44357
$ws.Range("J95").Value = 42
$ws.Range("K95").Value = 7000
$ws.Range("L95").Value = 8000
$ws.Range("M95").Value = 7452
$ws.Range("P95").Value = 2484
$ws.Range("D96").Value = 44215
$ws.Range("J96").Value = 160
$ws.Range("K96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("M96").Value = 10000
$ws.Range("P96").Value = 3333
$ws.Range("D97").Value = 44166
$ws.Range("J97").Value = 36
$ws.Range("K97").Value = 12000
$ws.Range("L97").Value = 12000
$ws.Range("M97").Value = 12000
$ws.Range("P97").Value = 4000
$ws.Range("D98").Value = 44259
$ws.Range("J98").Value = 39
$ws.Range("D99").Value = 44278
$ws.Range("J99").Value = 36
$ws.Range("D100").Value = 44218
$ws.Range("J100").Value = 38
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = 10000
$ws.Range("P100").Value = 3333
$ws.Range("D101").Value = 44280
$ws.Range("J101").Value = 37
$ws.Range("K101").Value = 10000
$ws.Range("L101").Value = 10000
$ws.Range("M101").Value = 10000
$ws.Range("P101").Value = 3333
$ws.Range("D102").Value = 44427
$ws.Range("J102").Value = 29
$ws.Range("K102").Value = 9000
$ws.Range("M102").Value = 9448
$ws.Range("P102").Value = 3149
$ws.Range("D103").Value = 44340
$ws.Range("J103").Value = 37
$ws.Range("K103").Value = 8000
$ws.Range("L103").Value = 9000
$ws.Range("M103").Value = 8405
$ws.Range("P103").Value = 2802
$ws.Range("D104").Value = 44497
$ws.Range("J104").Value = 36
$ws.Range("K104").Value = 8000
$ws.Range("L104").Value = 9000
$ws.Range("M104").Value = 8472
$ws.Range("P104").Value = 2824
$ws.Range("D105").Value = 44267
$ws.Range("J105").Value = 37
$ws.Range("K105").Value = 10000
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = 10000
$ws.Range("P105").Value = 3333
$ws.Range("D106").Value = 44418
$ws.Range("J106").Value = 29
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = 9448
$ws.Range("P106").Value = 3149
$ws.Range("D107").Value = 44525
$ws.Range("J107").Value = 35
$ws.Range("K107").Value = 8500
$ws.Range("L107").Value = 9000
$ws.Range("M107").Value = 8729
$ws.Range("P107").Value = 2910
$ws.Range("D108").Value = 44383
$ws.Range("J108").Value = 33
$ws.Range("K108").Value = 8500
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 8727
$ws.Range("P108").Value = 2909
$ws.Range("D109").Value = 44307
$ws.Range("J109").Value = 43
$ws.Range("K109").Value = 9000
$ws.Range("L109").Value = 10000
$ws.Range("M109").Value = 9581
$ws.Range("P109").Value = 3194
$ws.Range("D110").Value = 44299
$ws.Range("J110").Value = 38
$ws.Range("K110").Value = 9000
$ws.Range("M110").Value = 9579
$ws.Range("P110").Value = 3193
$ws.Range("D111").Value = 44316
$ws.Range("J111").Value = 36
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = 9000
$ws.Range("P111").Value = 3000
$ws.Range("D112").Value = 44229
$ws.Range("J112").Value = 33
$ws.Range("K112").Value = 10000
$ws.Range("M112").Value = 10000
$ws.Range("P112").Value = 3333
$ws.Range("D113").Value = 44239
$ws.Range("J113").Value = 28
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 10000
$ws.Range("P113").Value = 3333
$ws.Range("D114").Value = 44372
$ws.Range("J114").Value = 29
$ws.Range("M114").Value = 8414
$ws.Range("P114").Value = 2805
$ws.Range("D115").Value = 44469
$ws.Range("J115").Value = 32
$ws.Range("M115").Value = 8469
$ws.Range("P115").Value = 2823
$ws.Range("D116").Value = 44392
$ws.Range("J116").Value = 32
$ws.Range("M116").Value = 8469
$ws.Range("P116").Value = 2823
$ws.Range("D117").Value = 44273
$ws.Range("J117").Value = 39
$ws.Range("K117").Value = 10000
$ws.Range("L117").Value = 10000
$ws.Range("M117").Value = 10000
$ws.Range("P117").Value = 3333
$ws.Range("D118").Value = 44350
$ws.Range("J118").Value = 45
$ws.Range("K118").Value = 8000
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 8422
$ws.Range("P118").Value = 2807
$ws.Range("D119").Value = 44358
$ws.Range("J119").Value = 34
$ws.Range("K119").Value = 7000
$ws.Range("L119").Value = 8000
$ws.Range("M119").Value = 7441
$ws.Range("P119").Value = 2480
$ws.Range("D120").Value = 44189
$ws.Range("J120").Value = 44
$ws.Range("D121").Value = 44453
$ws.Range("J121").Value = 32
$ws.Range("M121").Value = 8438
$ws.Range("P121").Value = 2813
$ws.Range("D122").Value = 44314
$ws.Range("J122").Value = 22
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = 9000
$ws.Range("P122").Value = 3000
$ws.Range("D123").Value = 44509
$ws.Range("J123").Value = 34
$ws.Range("K123").Value = 8000
$ws.Range("L123").Value = 9000
$ws.Range("M123").Value = 8441
$ws.Range("P123").Value = 2814
$ws.Range("D124").Value = 44474
$ws.Range("J124").Value = 34
$ws.Range("K124").Value = 8500
$ws.Range("L124").Value = 9000
$ws.Range("M124").Value = 8735
$ws.Range("P124").Value = 2912
$ws.Range("D125").Value = 44224
$ws.Range("J125").Value = 38
$ws.Range("D126").Value = 44447
$ws.Range("J126").Value = 28
$ws.Range("M126").Value = 9464
$ws.Range("P126").Value = 3155
$ws.Range("D127").Value = 44274
$ws.Range("J127").Value = 25
$ws.Range("K127").Value = 10000
$ws.Range("L127").Value = 10000
$ws.Range("M127").Value = 10000
$ws.Range("P127").Value = 3333
$ws.Range("D128").Value = 44490
$ws.Range("J128").Value = 34
$ws.Range("K128").Value = 8500
$ws.Range("L128").Value = 9000
$ws.Range("M128").Value = 8735
$ws.Range("P128").Value = 2912
